$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$hf = $m.HeadersFooters
$dt = $hf.DateAndTime
Write-Host "Before: " $dt.Text
$dt.Text = "2021-08-06"
Write-Host "After direct: " $dt.Text
$hf2 = $m.HeadersFooters
$dt2 = $hf2.DateAndTime
Write-Host "After refetch: " $dt2.Text
